$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing input values (row 8-13) ---
$ws.Range("B8").Value = 470
$ws.Range("B9").Value = 470
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 10
$ws.Range("B13").Value = 150

# --- New header row 32: fg / kHz | R / Ohm | C / nF (grey fill) ---
# (populated before rows 27/28 so new shared-strings are appended in the
# same order the original authoring session produced them: R / Ohm, C / nF,
# then P / mW, U / V)
$ws.Range("A32").Value = "fg / kHz"
$ws.Range("B32").Value = "R / Ohm"
$ws.Range("C32").Value = "C / nF"
$ws.Range("A32:C32").Interior.Color = 14277081

# --- New rows 27/28: U / V and P / mW ---
$ws.Range("A28").Value = "P / mW"
$ws.Range("A27").Value = "U / V"
$ws.Range("B27").Value = 960
$ws.Range("B28").Formula = "=B27^2/SUM(B8:B11)"
$ws.Range("B28").NumberFormat = "0"

# --- Row 33: A,B given, C computed (0.000 format) ---
$ws.Range("A33").Value = 100
$ws.Range("B33").Value = 75

# --- Row 34: A,C given, B computed (0 format) ---
$ws.Range("A34").Value = 100
$ws.Range("B34").Formula = "=1/(2*PI()*A34*1000*C34*0.000000001)"
$ws.Range("B34").Interior.Color = 14277081
$ws.Range("B34").NumberFormat = "0"
$ws.Range("C34").Value = 15

# --- Row 35: B,C given, A computed (0 format) ---
$ws.Range("A35").Formula = "=1/(2*PI()*B35*C35*0.000000001)"
$ws.Range("A35").Interior.Color = 14277081
$ws.Range("A35").NumberFormat = "0"
$ws.Range("B35").Value = 75
$ws.Range("C35").Value = 18

# Row 33 C computed last, so the "0.000" numFmt is created after the "0" numFmt
# above (keeps generated style index order matching: s=5 integer, s=6 decimal)
$ws.Range("C33").Formula = "=1/(2*PI()*B33*A33*1000)*1000000000"
$ws.Range("C33").Interior.Color = 14277081
$ws.Range("C33").NumberFormat = "0.000"

# --- Picture (Grafik 1): move/resize to new anchor ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 314.0
$shp.Top = 227.0
$shp.Width = 423.5217322834646
$shp.Height = 308.01582677165356

# --- View state: scroll + selection ---
$ws.Range("D23").Select() | Out-Null
